$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert 4 blank rows before row 16 so the old "Resolution" block (rows
#    16-17) and the hyperlink row (28) shift down to rows 20-21 and 32,
#    leaving rows 16-17 free for the new "Lines" / "Bytes used" rows.
# ---------------------------------------------------------------------------
$ws.Rows("16:19").Insert()

# ---------------------------------------------------------------------------
# 2. Extend the table with a new column E ("1 MHz"). Copy formatting from
#    column D so the new column matches the look of the existing columns.
# ---------------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("D3:D14").Copy()
$ws.Range("E3:E14").PasteSpecial(-4122)

$ws.Range("E1").Value = "1 MHz"

# Row 10 used to contain the text "same" for columns C and D; it now holds
# the literal resolution value, repeated across C, D and the new E column.
$ws.Range("C10").Value = 480
$ws.Range("D10").Value = 480
$ws.Range("E10").Value = 480

# Formulas for the new "1 MHz" column.
$ws.Range("E3").Formula = "=B3/20"
$ws.Range("E4").Formula = "=B4/20"
$ws.Range("E5").Formula = "=B5/20"
$ws.Range("E6").Formula = "=B6/20"
$ws.Range("E7").Formula = "=B7/20"

# ---------------------------------------------------------------------------
# 3. New rows 16-17: "Lines" and "Bytes used (1 byte/pixel)".
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Lines"
$ws.Range("B16").Formula = "=B10"
$ws.Range("C16").Formula = "=C10/2"
$ws.Range("D16").Formula = "=`$D10/4"
$ws.Range("E16").Formula = "=`$D10/5"

$ws.Range("A17").Value = "Bytes used (1 byte/pixel)"
$ws.Range("B17").Formula = "=B3*B16"
$ws.Range("C17").Formula = "=C3*C16"
$ws.Range("D17").Formula = "=`$D3*D16"
$ws.Range("E17").Formula = "=`$D3*E16"

# ---------------------------------------------------------------------------
# 4. Header row (row 1) is centered - apply to the whole A1:E1 range, which
#    already carries the bold/filled/bordered "label" style copied above.
# ---------------------------------------------------------------------------
$ws.Range("A1:E1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Column A needs to be wide enough to fit the new labels.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 23.5703125

# ---------------------------------------------------------------------------
# 6. Fix up the hyperlink, which used to live in C28 and now lives in C32
#    after the row insert (the insert operation does not move the
#    hyperlink definition itself).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C32"), "http://lucidscience.com/pro-vga%20video%20generator-7.aspx")
$ws.Range("C32").Style = "Hiperlink"

# ---------------------------------------------------------------------------
# 7. Selection / window bookkeeping to match the saved state.
# ---------------------------------------------------------------------------
$ws.Range("F14").Select()

$win = $excel.ActiveWindow
$win.Top = 165
$win.Height = 11250
